$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Collection_QRS_EQ5D-5L")

# Insert a new column at K (shifting short_name .. change_history one column
# to the right) and give the new column its header, "categories".
$jWidth = $ws.Columns("J:J").ColumnWidth
$ws.Columns("K:K").Insert() | Out-Null
$ws.Columns("K:K").ColumnWidth = $jWidth
$ws.Range("K1").Value = "categories"

# Re-apply the auto filter so the range now covers the new column (A1:AH16)
$ws.AutoFilterMode = $false
$null = $ws.Range("A1:AH16").AutoFilter()

# Keep the hidden _FilterDatabase defined name in sync with the wider range
foreach ($n in $wb.Names) {
    if ($n.Name -like "*_FilterDatabase*") {
        $n.RefersTo = "='Collection_QRS_EQ5D-5L'!`$A`$1:`$AH`$16"
    }
}

$null = $ws.Range("K6").Select()
